$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 403.48648
$ws.Range("I33").Value = 223.46875
$ws.Range("J33").Value = 1555.6
$ws.Range("K33").Value = 223.46875
$ws.Range("L33").Value = 1555.6
$ws.Range("M33").Value = 5.53125
$ws.Range("N33").Value = -2013.6

# Row 58
$ws.Range("H58").Value = 900
$ws.Range("I58").Value = 812.5
$ws.Range("J58").Value = 1133.3334
$ws.Range("K58").Value = 2437.5
$ws.Range("L58").Value = 3400.0002
$ws.Range("M58").Value = -2287.5
$ws.Range("N58").Value = -3700.0002

# Row 76
$ws.Range("H76").Value = 3036.3635
$ws.Range("I76").Value = 3022.2222
$ws.Range("K76").Value = 3022.2222
$ws.Range("M76").Value = -2707.2222

# Row 79
$ws.Range("H79").Value = 3036.3635
$ws.Range("I79").Value = 3022.2222
$ws.Range("K79").Value = 3022.2222
$ws.Range("M79").Value = -1930.2222

# Row 98
$ws.Range("H98").Value = 1532.3667
$ws.Range("I98").Value = 1474.36
$ws.Range("J98").Value = 1822.4
$ws.Range("K98").Value = 1474.36
$ws.Range("L98").Value = 1822.4
$ws.Range("M98").Value = 23.6400000000001
$ws.Range("N98").Value = -4818.4

# Row 116
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 442
$ws.Range("N116").ClearContents()

# Row 122
$ws.Range("H122").Value = 1532.3667
$ws.Range("I122").Value = 1474.36
$ws.Range("J122").Value = 1822.4
$ws.Range("K122").Value = 4423.08
$ws.Range("L122").Value = 5467.200000000001
$ws.Range("M122").Value = -1973.08
$ws.Range("N122").Value = -10367.2

# Row 132
$ws.Range("H132").Value = 4378.675
$ws.Range("I132").Value = 1296.0883
$ws.Range("J132").Value = 21846.666
$ws.Range("K132").Value = 3888.2649
$ws.Range("L132").Value = 65539.99800000001
$ws.Range("M132").Value = -1358.2649
$ws.Range("N132").Value = -70599.99800000001

# Row 138
$ws.Range("H138").Value = 6932862
$ws.Range("I138").Value = 1622.84
$ws.Range("J138").Value = 50253108
$ws.Range("K138").Value = 4868.52
$ws.Range("L138").Value = 150759324
$ws.Range("M138").Value = 271.4800000000005
$ws.Range("N138").Value = -150769604

# Row 141
$ws.Range("H141").Value = 3484.3333
$ws.Range("I141").Value = 1979.1111
$ws.Range("K141").Value = 5937.3333
$ws.Range("M141").Value = -757.3333000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 1533
$ws.Range("I3").Value = 1533
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1533
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1418
$ws.Range("N3").ClearContents()

# Row 32
$ws.Range("H32").Value = 5569
$ws.Range("I32").Value = 3164.4219
$ws.Range("J32").Value = 17406.924
$ws.Range("K32").Value = 3164.4219
$ws.Range("L32").Value = 17406.924
$ws.Range("M32").Value = -2877.4219
$ws.Range("N32").Value = -17980.924

# Row 63
$ws.Range("H63").Value = 3438.389
$ws.Range("I63").Value = 2209.4443
$ws.Range("J63").Value = 4667.3335
$ws.Range("K63").Value = 2209.4443
$ws.Range("L63").Value = 4667.3335
$ws.Range("M63").Value = -1523.4443
$ws.Range("N63").Value = -6039.3335

# Row 66
$ws.Range("H66").Value = 3438.389
$ws.Range("I66").Value = 2209.4443
$ws.Range("J66").Value = 4667.3335
$ws.Range("K66").Value = 11047.2215
$ws.Range("L66").Value = 23336.6675
$ws.Range("M66").Value = -7615.2215
$ws.Range("N66").Value = -30200.6675

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1964.4138
$ws.Range("I99").Value = 1504.25
$ws.Range("J99").Value = 2530.7693
$ws.Range("K99").Value = 1504.25
$ws.Range("L99").Value = 2530.7693
$ws.Range("M99").Value = -6.25
$ws.Range("N99").Value = -5526.7693

# Row 105
$ws.Range("H105").Value = 2389.8838
$ws.Range("I105").Value = 2384.5122
$ws.Range("K105").Value = 2384.5122
$ws.Range("M105").Value = -637.5122000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1782.375
$ws.Range("I31").Value = 1341.4615
$ws.Range("K31").Value = 1341.4615
$ws.Range("M31").Value = -1046.4615

# Row 34
$ws.Range("H34").Value = 1782.375
$ws.Range("I34").Value = 1341.4615
$ws.Range("K34").Value = 1341.4615
$ws.Range("M34").Value = -1139.4615

# Row 99
$ws.Range("H99").Value = 2387.4075
$ws.Range("I99").Value = 2157.2727
$ws.Range("J99").Value = 3400
$ws.Range("K99").Value = 2157.2727
$ws.Range("L99").Value = 3400
$ws.Range("M99").Value = -659.2727
$ws.Range("N99").Value = -6396

# Row 103
$ws.Range("H103").Value = 25789.143
$ws.Range("I103").Value = 5524
$ws.Range("J103").Value = 29166.666
$ws.Range("K103").Value = 5524
$ws.Range("L103").Value = 29166.666
$ws.Range("M103").Value = -4352
$ws.Range("N103").Value = -31510.666

# Row 126
$ws.Range("H126").Value = 2387.4075
$ws.Range("I126").Value = 2157.2727
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 6471.8181
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -4001.8181
$ws.Range("N126").Value = -15140

# Row 134
$ws.Range("H134").Value = 2902.0715
$ws.Range("I134").Value = 1302.7273
$ws.Range("K134").Value = 3908.1819
$ws.Range("M134").Value = -1373.1819

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 17080.678
$ws.Range("I4").Value = 220.2
$ws.Range("J4").Value = 20323.076
$ws.Range("K4").Value = 660.5999999999999
$ws.Range("L4").Value = 60969.228
$ws.Range("M4").Value = -548.5999999999999
$ws.Range("N4").Value = -61193.228

# Row 131
$ws.Range("H131").Value = 925.60205
$ws.Range("J131").Value = 940.71875
$ws.Range("L131").Value = 2822.15625
$ws.Range("N131").Value = -12902.15625

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 2978
$ws.Range("I4").Value = 1900
$ws.Range("J4").Value = 3247.5
$ws.Range("K4").Value = 1900
$ws.Range("L4").Value = 3247.5
$ws.Range("M4").Value = -1788
$ws.Range("N4").Value = -3471.5

# Row 5
$ws.Range("H5").Value = 8080.5713
$ws.Range("J5").Value = 9166
$ws.Range("L5").Value = 9166
$ws.Range("N5").Value = -9390

# Row 57
$ws.Range("H57").Value = 19244.223
$ws.Range("J57").Value = 19244.223
$ws.Range("L57").Value = 19244.223
$ws.Range("N57").Value = -20884.223

# Row 70
$ws.Range("H70").Value = 5335.35
$ws.Range("I70").Value = 5250.5713
$ws.Range("J70").Value = 5533.1665
$ws.Range("K70").Value = 5250.5713
$ws.Range("L70").Value = 5533.1665
$ws.Range("M70").Value = -4980.5713
$ws.Range("N70").Value = -6073.1665

# Row 73
$ws.Range("H73").Value = 5335.35
$ws.Range("I73").Value = 5250.5713
$ws.Range("J73").Value = 5533.1665
$ws.Range("K73").Value = 5250.5713
$ws.Range("L73").Value = 5533.1665
$ws.Range("M73").Value = -4314.5713
$ws.Range("N73").Value = -7405.1665

# Row 124
$ws.Range("H124").Value = 27789.475
$ws.Range("J124").Value = 27789.475
$ws.Range("L124").Value = 27789.475
$ws.Range("N124").Value = -37609.475

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1256.04
$ws.Range("I7").Value = 1221.75
$ws.Range("J7").Value = 1317
$ws.Range("K7").Value = 1221.75
$ws.Range("L7").Value = 1317
$ws.Range("M7").Value = -1109.75
$ws.Range("N7").Value = -1541

# Row 126
$ws.Range("H126").Value = 1256.04
$ws.Range("I126").Value = 1221.75
$ws.Range("J126").Value = 1317
$ws.Range("K126").Value = 3665.25
$ws.Range("L126").Value = 3951
$ws.Range("M126").Value = -1195.25
$ws.Range("N126").Value = -8891

# Row 132
$ws.Range("H132").Value = 4816.353
$ws.Range("I132").Value = 4988.2
$ws.Range("J132").Value = 4570.857
$ws.Range("K132").Value = 14964.6
$ws.Range("L132").Value = 13712.571
$ws.Range("M132").Value = -12434.6
$ws.Range("N132").Value = -18772.571

# Row 136
$ws.Range("H136").Value = 6848.8423
$ws.Range("I136").Value = 1329.1428
$ws.Range("J136").Value = 22304
$ws.Range("K136").Value = 3987.4284
$ws.Range("L136").Value = 66912
$ws.Range("M136").Value = -1437.4284
$ws.Range("N136").Value = -72012

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 386666.66
$ws.Range("I4").Value = 80000
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 80000
$ws.Range("L4").Value = 1000000
$ws.Range("M4").Value = -79887
$ws.Range("N4").Value = -1000226

